$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Numeric value updates in column E (bound/limit values) ---
$ws.Range("E20").Value = 40
$ws.Range("E22").Value = 20
$ws.Range("E24").Value = 0
$ws.Range("E29").Value = 40
$ws.Range("E31").Value = 15
$ws.Range("E33").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("E35").Value = 3.3
$ws.Range("E36").Value = 4.0999999999999996
$ws.Range("E37").Value = 4.0999999999999996
$ws.Range("E38").Value = 4.0999999999999996
$ws.Range("E40").Value = 5
$ws.Range("E42").Value = 5
$ws.Range("E43").Value = 7
$ws.Range("E44").Value = 8
$ws.Range("E45").Value = 9
$ws.Range("E47").Value = 12
$ws.Range("E49").Value = 15
$ws.Range("E60").Value = 7
$ws.Range("E61").Value = 15
$ws.Range("E62").Value = 25
$ws.Range("E63").Value = 40
$ws.Range("E64").Value = 50
$ws.Range("E65").Value = 58
$ws.Range("E67").Value = 60
$ws.Range("E69").Value = 14
$ws.Range("E70").Value = 20
$ws.Range("E71").Value = 23
$ws.Range("E72").Value = 28
$ws.Range("E73").Value = 32
$ws.Range("E74").Value = 35
$ws.Range("E76").Value = 35

# --- Installed capacity style rows (previously all 1000) ---
$ws.Range("E78").Value = 0
$ws.Range("E79").Value = 10
$ws.Range("E80").Value = 20
$ws.Range("E81").Value = 40
$ws.Range("E82").Value = 50
$ws.Range("E83").Value = 60
$ws.Range("E84").Value = 80

$ws.Range("E86").Value = 0
$ws.Range("E87").Value = 8
$ws.Range("E88").Value = 20
$ws.Range("E89").Value = 50
$ws.Range("E90").Value = 100
$ws.Range("E91").Value = 120
$ws.Range("E92").Value = 150

# --- Fuel price forecast values ---
$ws.Range("E94").Value = 141.68
$ws.Range("E95").Value = 123.2
$ws.Range("E96").Value = 115.5
$ws.Range("E97").Value = 107.8
$ws.Range("E98").Value = 101.64
$ws.Range("E99").Value = 96.25

# --- Row 100 header: add attribute/unit labels (matches row 93 pattern) ---
$ws.Range("C100").Value = "Prognoza ceny paliwa"
$ws.Range("E100").Value = "zł/GJ"

$ws.Range("E101").Value = 765.45
$ws.Range("E102").Value = 800.87
$ws.Range("E103").Value = 840.91
$ws.Range("E104").Value = 881.29
$ws.Range("E105").Value = 921.96
$ws.Range("E106").Value = 952.4

# --- Clear the highlight fills that previously marked these cells ---
# (orange FFC000 and green 92D050 highlighted ranges become plain white)
# NOTE: the COM engine only honours the first area of a multi-area Range,
# so each contiguous block is set individually.
$ws.Range("E15:E20").Interior.Color = 16777215
$ws.Range("E22").Interior.Color = 16777215
$ws.Range("E24:E29").Interior.Color = 16777215
$ws.Range("E31").Interior.Color = 16777215
$ws.Range("E33:E38").Interior.Color = 16777215
$ws.Range("E40").Interior.Color = 16777215
$ws.Range("E42:E47").Interior.Color = 16777215
$ws.Range("E49").Interior.Color = 16777215
$ws.Range("E51:E56").Interior.Color = 16777215
$ws.Range("E58").Interior.Color = 16777215
$ws.Range("E60:E65").Interior.Color = 16777215
$ws.Range("E67").Interior.Color = 16777215
$ws.Range("E69:E74").Interior.Color = 16777215
$ws.Range("E78:E84").Interior.Color = 16777215
$ws.Range("E86:E92").Interior.Color = 16777215
$ws.Range("E94:E99").Interior.Color = 16777215
$ws.Range("E101:E106").Interior.Color = 16777215
